$d = $word.ActiveDocument

# Locate the misspelled word "Bodyshiedet" (it occurs once in the document;
# a correctly spelled "Bodyshieldet" already exists later in the same
# paragraph). The fix inserts the missing "l" and, per the target
# revision, the resulting word is split across three runs:
#   "Bodyshie" + "l" + "det"
$rng = $d.Content
$found = $rng.Find.Execute("Bodyshiedet", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $matchStart = $rng.Start
    $matchEnd = $rng.End

    # Leave the first 8 characters ("Bodyshie") exactly as they were -- this
    # keeps the run (and the surrounding spellStart/spellEnd proofErr pair)
    # anchored in place. Replace only the trailing "det" with two new runs
    # ("l" and "det"), which is enough to reproduce the three-run split
    # described by the edit while keeping proofErr correctly wrapping the
    # whole word.
    $tail = $d.Range($matchStart + 8, $matchEnd)

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>l</w:t></w:r><w:r><w:t>det</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $tail.InsertXML($xml)
}
